$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2025-09-01 01:14:00"
$ws.Range("E3").Value = "2025-09-01 01:14:43"
$ws.Range("E4").Value = "2025-09-01 01:15:27"
$ws.Range("E5").Value = "2025-09-01 01:16:10"
$ws.Range("E6").Value = "2025-09-01 01:16:54"
